$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 62
$ws.Range("H62").Value = 14038650
$ws.Range("I62").Value = 20516350
$ws.Range("K62").Value = 20516350
$ws.Range("M62").Value = -20515726

# row 65
$ws.Range("H65").Value = 14038650
$ws.Range("I65").Value = 20516350
$ws.Range("K65").Value = 102581750
$ws.Range("M65").Value = -102578630

# row 70
$ws.Range("H70").Value = 1371762.4
$ws.Range("I70").Value = 4362.5
$ws.Range("J70").Value = 2587229
$ws.Range("K70").Value = 13087.5
$ws.Range("L70").Value = 7761687
$ws.Range("M70").Value = -12817.5
$ws.Range("N70").Value = -7762227

# row 73
$ws.Range("H73").Value = 1371762.4
$ws.Range("I73").Value = 4362.5
$ws.Range("J73").Value = 2587229
$ws.Range("K73").Value = 13087.5
$ws.Range("L73").Value = 7761687
$ws.Range("M73").Value = -12151.5
$ws.Range("N73").Value = -7763559

# row 106
$ws.Range("H106").Value = 3181.4546
$ws.Range("I106").Value = 3013.7144
$ws.Range("K106").Value = 3013.7144
$ws.Range("M106").Value = -2382.7144

# row 113
$ws.Range("H113").Value = 2599.4
$ws.Range("I113").Value = 2466
$ws.Range("J113").Value = 2799.5
$ws.Range("K113").Value = 2466
$ws.Range("L113").Value = 2799.5
$ws.Range("M113").Value = 788
$ws.Range("N113").Value = -9307.5

# row 132
$ws.Range("H132").Value = 1272.381
$ws.Range("I132").Value = 1080
$ws.Range("K132").Value = 3240
$ws.Range("M132").Value = -710

# row 138
$ws.Range("H138").Value = 2799
$ws.Range("I138").Value = 2183.6
$ws.Range("J138").Value = 3140.889
$ws.Range("K138").Value = 6550.799999999999
$ws.Range("L138").Value = 9422.667000000001
$ws.Range("M138").Value = -1410.799999999999
$ws.Range("N138").Value = -19702.667

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 21852.7
$ws.Range("I32").Value = 21610.053
$ws.Range("K32").Value = 21610.053
$ws.Range("M32").Value = -21323.053

# row 61
$ws.Range("H61").Value = 9914.615
$ws.Range("I61").Value = 8094.8823
$ws.Range("J61").Value = 13351.889
$ws.Range("K61").Value = 8094.8823
$ws.Range("L61").Value = 13351.889
$ws.Range("M61").Value = -7882.8823
$ws.Range("N61").Value = -13775.889

# row 74
$ws.Range("H74").Value = 3225.7368
$ws.Range("I74").Value = 1066.125
$ws.Range("K74").Value = 1066.125
$ws.Range("M74").Value = -192.125

# row 77
$ws.Range("H77").Value = 3225.7368
$ws.Range("I77").Value = 1066.125
$ws.Range("K77").Value = 5330.625
$ws.Range("M77").Value = -962.625

# row 132
$ws.Range("H132").Value = 6157.6665
$ws.Range("I132").Value = 4560.231
$ws.Range("K132").Value = 13680.693
$ws.Range("M132").Value = -11150.693

# row 136
$ws.Range("H136").Value = 9914.615
$ws.Range("I136").Value = 8094.8823
$ws.Range("J136").Value = 13351.889
$ws.Range("K136").Value = 24284.6469
$ws.Range("L136").Value = 40055.667
$ws.Range("M136").Value = -21734.6469
$ws.Range("N136").Value = -45155.667

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 2724
$ws.Range("I22").Value = 1827.4286
$ws.Range("K22").Value = 1827.4286
$ws.Range("M22").Value = -1654.4286

# row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

# row 132
$ws.Range("H132").Value = 110663.5
$ws.Range("J132").Value = 110663.5
$ws.Range("L132").Value = 110663.5
$ws.Range("N132").Value = -120783.5

$ws = $wb.Worksheets.Item("CRP")
# row 58
$ws.Range("H58").Value = 5622.737
$ws.Range("I58").Value = 4772.25
$ws.Range("J58").Value = 6241.273
$ws.Range("K58").Value = 4772.25
$ws.Range("L58").Value = 6241.273
$ws.Range("M58").Value = -4569.25
$ws.Range("N58").Value = -6647.273

# row 136
$ws.Range("H136").Value = 5622.737
$ws.Range("I136").Value = 4772.25
$ws.Range("J136").Value = 6241.273
$ws.Range("K136").Value = 14316.75
$ws.Range("L136").Value = 18723.819
$ws.Range("M136").Value = -11766.75
$ws.Range("N136").Value = -23823.819

$ws = $wb.Worksheets.Item("CUL")
# row 131
$ws.Range("H131").Value = 8774976
$ws.Range("J131").Value = 3128.3784
$ws.Range("L131").Value = 9385.135200000001
$ws.Range("N131").Value = -19465.1352

# row 132
$ws.Range("H132").Value = 68782.266
$ws.Range("I132").Value = 126125.125
$ws.Range("J132").Value = 3247.5715
$ws.Range("K132").Value = 1135126.125
$ws.Range("L132").Value = 29228.1435
$ws.Range("M132").Value = -1132596.125
$ws.Range("N132").Value = -34288.1435

# row 139
$ws.Range("H139").Value = 2817.5454
$ws.Range("I139").Value = 2494.2856
$ws.Range("K139").Value = 7482.8568
$ws.Range("M139").Value = -2342.8568

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 4999.8
$ws.Range("I70").Value = 4999.8
$ws.Range("K70").Value = 4999.8
$ws.Range("M70").Value = -4729.8

# row 73
$ws.Range("H73").Value = 4999.8
$ws.Range("I73").Value = 4999.8
$ws.Range("K73").Value = 4999.8
$ws.Range("M73").Value = -4063.8

# row 97
$ws.Range("H97").Value = 1110.7142
$ws.Range("I97").Value = 1050
$ws.Range("K97").Value = 1050
$ws.Range("M97").Value = -554

# row 122
$ws.Range("H122").Value = 1418.8
$ws.Range("I122").Value = 1593
$ws.Range("K122").Value = 4779
$ws.Range("M122").Value = -2329

# row 126
$ws.Range("H126").Value = 3911.0908
$ws.Range("I126").Value = 2181
$ws.Range("K126").Value = 6543
$ws.Range("M126").Value = -4073

# row 132
$ws.Range("H132").Value = 6387.421
$ws.Range("I132").Value = 3223.7273
$ws.Range("K132").Value = 9671.1819
$ws.Range("M132").Value = -7141.1819

# row 136
$ws.Range("H136").Value = 72945.28999999999
$ws.Range("J136").Value = 72945.28999999999
$ws.Range("L136").Value = 218835.87
$ws.Range("N136").Value = -223935.87

$ws = $wb.Worksheets.Item("LTW")
# row 132
$ws.Range("H132").Value = 12170.75
$ws.Range("I132").Value = 8100
$ws.Range("J132").Value = 16241.5
$ws.Range("K132").Value = 24300
$ws.Range("L132").Value = 48724.5
$ws.Range("M132").Value = -21770
$ws.Range("N132").Value = -53784.5

# row 136
$ws.Range("H136").Value = 4765.271
$ws.Range("I136").Value = 4233.7715
$ws.Range("K136").Value = 12701.3145
$ws.Range("M136").Value = -10151.3145

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 6003
$ws.Range("J62").Value = 6003
$ws.Range("L62").Value = 6003
$ws.Range("N62").Value = -7251

# row 65
$ws.Range("H65").Value = 6003
$ws.Range("J65").Value = 6003
$ws.Range("L65").Value = 30015
$ws.Range("N65").Value = -36255

# row 81
$ws.Range("H81").Value = 4426.846
$ws.Range("I81").Value = 4123.96
$ws.Range("K81").Value = 8247.92
$ws.Range("M81").Value = -7186.92

# row 84
$ws.Range("H84").Value = 4426.846
$ws.Range("I84").Value = 4123.96
$ws.Range("K84").Value = 41239.6
$ws.Range("M84").Value = -35935.6

# row 126
$ws.Range("H126").Value = 38144.832
$ws.Range("I126").Value = 50098.684
$ws.Range("K126").Value = 150296.052
$ws.Range("M126").Value = -147826.052

# row 132
$ws.Range("H132").Value = 7052.8823
$ws.Range("I132").Value = 5836.25
$ws.Range("J132").Value = 8134.3335
$ws.Range("K132").Value = 17508.75
$ws.Range("L132").Value = 24403.0005
$ws.Range("M132").Value = -14978.75
$ws.Range("N132").Value = -29463.0005

# row 136
$ws.Range("H136").Value = 3805.9167
$ws.Range("I136").Value = 1512.2
$ws.Range("J136").Value = 7628.778
$ws.Range("K136").Value = 4536.6
$ws.Range("L136").Value = 22886.334
$ws.Range("M136").Value = -1986.6
$ws.Range("N136").Value = -27986.334
